$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.015000000000001
$ws.Range("A3").Value = -21.603
$ws.Range("B5").Value = 6.078
$ws.Range("E5").Value = 13.132
$ws.Range("E9").Value = 13.238
$ws.Range("E11").Value = 13.132
$ws.Range("A14").Value = -20.719
$ws.Range("A16").Value = -20.905
$ws.Range("B16").Value = 6.566
$ws.Range("E17").Value = 13.644
$ws.Range("A21").Value = -20.993
$ws.Range("E21").Value = 13.357
$ws.Range("A23").Value = -21.584
$ws.Range("A25").Value = -22.27
